# Adds a new "targetThicknesses" worksheet (with a small table of target
# thicknesses / implanted-atom counts) at the end of the workbook, and
# makes it the active/selected sheet - matching the commit
# "Added tables and a bit more explanation to target characterization section."

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("implantedTargets")
# so it lands at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "targetThicknesses"

# Column A: target labels
$newSheet.Range("A1").Value = "Target"
$newSheet.Range("A2").Value = "TiN"
$newSheet.Range("A3").Value = "ZrN 5"
$newSheet.Range("A4").Value = "ZrN 1"
$newSheet.Range("A5").Value = "ZrN 12"

# Column B: implanted-atom areal densities
$newSheet.Range("B1").Value = "`$n`$ \left(10^{17} atoms / cm^{2} \right)"
$newSheet.Range("B2").Value = "7.070 `$\pm`$ 0.566"
$newSheet.Range("B3").Value = "5.623 `$\pm`$ 0.450"
$newSheet.Range("B4").Value = "5.826 `$\pm`$ 0.466"
$newSheet.Range("B5").Value = "10.480 `$\pm`$ 0.834"

# Make the new sheet the active tab with the whole table selected.
$newSheet.Select()
$newSheet.Range("A1:B5").Select() | Out-Null
